# UI tweaks and data update
# - Update surveyed coordinate rows 111-116 with corrected readings
# - Append newly surveyed rows 117-123
# - Update the view: scroll position / selected cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column letter -> index map used below: A=1 B=2 C=3 E=5 F=6 (D/G/H/I stay blank)
$rowData = @(
    @(111, @(1, 1818764.4280000001, 2, 5578108.6739999996, 3, 287.20499999999998, 5, 1.123, 6, 48.885000000000005)),
    @(112, @(1, 1818764.024, 2, 5578109.9699999997, 3, 286.97800000000001, 5, 1.3580000000000001, 6, 50.243000000000002)),
    @(113, @(1, 1818763.557, 2, 5578111.71, 3, 286.57400000000001, 5, 1.802, 6, 52.045000000000002)),
    @(114, @(1, 1818762.865, 2, 5578113.3329999996, 3, 286.03699999999998, 5, 1.764, 6, 53.809000000000005)),
    @(115, @(1, 1818762.585, 2, 5578114.5259999996, 3, 285.64499999999998, 5, 1.2250000000000001, 6, 55.034000000000006)),
    @(116, @(1, 1818762.3019999999, 2, 5578116.159, 3, 285.61799999999999, 5, 1.657, 6, 56.691000000000003)),
    @(117, @(1, 1818761.8430000001, 2, 5578118.1330000004, 3, 285.58800000000002, 5, 2.0270000000000001, 6, 58.718000000000004)),
    @(118, @(1, 1818761.379, 2, 5578119.4029999999, 3, 285.59199999999998, 5, 1.3520000000000001, 6, 60.07)),
    @(119, @(1, 1818761.3370000001, 2, 5578119.8150000004, 3, 285.875, 5, 0.41399999999999998, 6, 60.484000000000002)),
    @(120, @(1, 1818761.2039999999, 2, 5578120.307, 3, 286.38600000000002, 5, 0.51, 6, 60.994)),
    @(121, @(1, 1818761.3430000001, 2, 5578120.665, 3, 287.101, 5, 0.38400000000000001, 6, 61.378)),
    @(122, @(1, 1818760.024, 2, 5578121.9440000001, 3, 287.41500000000002, 5, 1.837, 6, 63.215000000000003)),
    @(123, @(1, 1818758.6529999999, 2, 5578124.7249999996, 3, 288.31299999999999, 5, 3.101, 6, 66.316000000000003))
)

foreach ($item in $rowData) {
    $r = $item[0]
    $cells = $item[1]
    for ($i = 0; $i -lt $cells.Length; $i += 2) {
        $col = $cells[$i]
        $val = $cells[$i + 1]
        $ws.Cells.Item($r, $col).Value = $val
    }
}

# Scroll the view so row 94 is at the top, then select I119 (new active cell)
$excel.Goto($ws.Range("A94"), $true)
$ws.Range("I119").Select()
